$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (old N,O,P shift right to O,P,Q)
$ws.Columns("N:N").Insert()

# The new column inherits the width of its left neighbour (column M),
# losing the "best fit" flag but keeping an explicit custom width.
$ws.Columns("N:N").ColumnWidth = 9.8

# Re-assert the numeric values that were shifted into column Q so they keep
# their original (non drifted) two-decimal precision after the shift.
$ws.Range("Q4").Value = 887.72
$ws.Range("Q5").Value = 887.72
$ws.Range("Q6").Value = 887.72
$ws.Range("Q7").Value = 887.72
$ws.Range("Q8").Value = 715.92

# Update the active cell/selection as recorded after the edit
$ws.Range("K15").Select()
